$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume (E) columns are plain text in the source
# workbook (scraped values such as thousands-dot-separated prices and
# padded "  +x.xx%  " percentages). Volume-column values always stay
# text because of the leading/trailing spaces and "%" sign, and Price
# values written with two "." separators (e.g. "62.837.87") are never
# parsed as numbers either, so both can be assigned directly.
#
# Price values that look like an ordinary decimal number (one "." or
# none) would otherwise be auto-converted by Excel into a genuine
# number -- silently dropping a significant trailing zero in several
# cases (e.g. "70.50" -> 70.5). For those cells we briefly mark the
# cell as Text before assigning the value, then restore the default
# "Normal" style so the cell carries no visible/semantic formatting
# change -- only its literal text content differs, as in the source diff.

$ws.Range("D2").Value = "62.837.87"
$ws.Range("E2").Value = "  +4.79%  "
$ws.Range("D3").Value = "3.352.84"
$ws.Range("E3").Value = "  +4.99%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.530"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.65%  "
$ws.Range("E9").Value = "  +2.41%  "
$ws.Range("E10").Value = "  +4.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.439"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.86%  "
$ws.Range("D12").Value = "3.927.96"
$ws.Range("E12").Value = "  +4.81%  "
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.06%  "
$ws.Range("E15").Value = "  +3.64%  "
$ws.Range("D16").Value = "62.903.50"
$ws.Range("E16").Value = "  +4.77%  "
$ws.Range("D17").Value = "3.376.77"
$ws.Range("E17").Value = "  +5.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "387.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("E22").Value = "  +2.37%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("E25").Value = "  +4.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("E27").Value = "  +7.71%  "
$ws.Range("E29").Value = "  +7.40%  "
$ws.Range("E30").Value = "  +4.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.95%  "
$ws.Range("E33").Value = "  +7.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.70"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "160.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.90%  "
$ws.Range("E36").Value = "  +9.37%  "
$ws.Range("E37").Value = "  +12.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0741"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.74%  "
$ws.Range("D40").Value = "2.826.62"
$ws.Range("E40").Value = "  +1.80%  "
$ws.Range("E41").Value = "  +8.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.746"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.35%  "
$ws.Range("E45").Value = "  +4.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.04%  "
$ws.Range("D47").Value = "3.395.41"
$ws.Range("E47").Value = "  +4.93%  "
$ws.Range("E48").Value = "  +2.23%  "
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.807"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "280.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.19%  "
